$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9549537301063538
$ws.Range("B1").Value = 1.15474259853363
$ws.Range("C1").Value = 0.9737921953201294
$ws.Range("D1").Value = 0.9063471555709839
$ws.Range("E1").Value = 0.9425735473632812
